$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.16209533333333
$ws.Range("H2").Value = 51.486286
$ws.Range("I2").Value = 0.2459970657298922
$ws.Range("J2").Value = 0.2459970657298922
$ws.Range("M2").Value = 2.883158333333334
$ws.Range("N2").Value = 8.649475
$ws.Range("O2").Value = 0.1005826776766536
$ws.Range("P2").Value = 0.1005826776766536
$ws.Range("Q2").Value = 49.48103817776111
$ws.Range("R2").Value = 445.32934359985
$ws.Range("S2").Value = 0.02474304357171231
$ws.Range("T2").Value = 0.02474304357171231

$ws.Range("G3").Value = 17.16209533333333
$ws.Range("H3").Value = 51.486286
$ws.Range("I3").Value = 0.2459970657298922
$ws.Range("J3").Value = 0.2459970657298922
$ws.Range("O3").Value = 0.4854237085598054
$ws.Range("P3").Value = 0.4854237085598054
$ws.Range("Q3").Value = 238.8012489869643
$ws.Range("R3").Value = 2149.211240882678
$ws.Range("S3").Value = 0.1194128079414345
$ws.Range("T3").Value = 0.1194128079414345

$ws.Range("G4").Value = 17.16209533333333
$ws.Range("H4").Value = 51.486286
$ws.Range("I4").Value = 0.2459970657298922
$ws.Range("J4").Value = 0.2459970657298922
$ws.Range("M4").Value = 11.86694533333333
$ws.Range("N4").Value = 35.600836
$ws.Range("O4").Value = 0.413993613763541
$ws.Range("P4").Value = 0.413993613763541
$ws.Range("Q4").Value = 203.6616471261218
$ws.Range("R4").Value = 1832.954824135096
$ws.Range("S4").Value = 0.1018412142167454
$ws.Range("T4").Value = 0.1018412142167454

$ws.Range("I5").Value = 0.2105756965403629
$ws.Range("J5").Value = 0.2105756965403628
$ws.Range("M5").Value = 2.883158333333334
$ws.Range("N5").Value = 8.649475
$ws.Range("O5").Value = 0.1005826776766536
$ws.Range("P5").Value = 0.1005826776766536
$ws.Range("Q5").Value = 42.35621286338056
$ws.Range("R5").Value = 381.205915770425
$ws.Range("S5").Value = 0.02118026741165613
$ws.Range("T5").Value = 0.02118026741165613

$ws.Range("I6").Value = 0.2105756965403629
$ws.Range("J6").Value = 0.2105756965403628
$ws.Range("O6").Value = 0.4854237085598054
$ws.Range("P6").Value = 0.4854237085598054
$ws.Range("Q6").Value = 204.4160128127421
$ws.Range("S6").Value = 0.1022184355471871
$ws.Range("T6").Value = 0.1022184355471871

$ws.Range("I7").Value = 0.2105756965403629
$ws.Range("J7").Value = 0.2105756965403628
$ws.Range("M7").Value = 11.86694533333333
$ws.Range("N7").Value = 35.600836
$ws.Range("O7").Value = 0.413993613763541
$ws.Range("P7").Value = 0.413993613763541
$ws.Range("Q7").Value = 174.3361981773809
$ws.Range("R7").Value = 1569.025783596428
$ws.Range("S7").Value = 0.08717699358151962
$ws.Range("T7").Value = 0.08717699358151959

$ws.Range("G8").Value = 37.91244433333333
$ws.Range("H8").Value = 113.737333
$ws.Range("I8").Value = 0.543427237729745
$ws.Range("J8").Value = 0.543427237729745
$ws.Range("M8").Value = 2.883158333333334
$ws.Range("N8").Value = 8.649475
$ws.Range("O8").Value = 0.1005826776766536
$ws.Range("P8").Value = 0.1005826776766536
$ws.Range("Q8").Value = 109.3075798166861
$ws.Range("R8").Value = 983.7682183501752
$ws.Range("S8").Value = 0.05465936669328513
$ws.Range("T8").Value = 0.05465936669328513

$ws.Range("G9").Value = 37.91244433333333
$ws.Range("H9").Value = 113.737333
$ws.Range("I9").Value = 0.543427237729745
$ws.Range("J9").Value = 0.543427237729745
$ws.Range("O9").Value = 0.4854237085598054
$ws.Range("P9").Value = 0.4854237085598054
$ws.Range("Q9").Value = 527.5311017160233
$ws.Range("R9").Value = 4747.77991544421
$ws.Range("S9").Value = 0.2637924650711838
$ws.Range("T9").Value = 0.2637924650711838

$ws.Range("G10").Value = 37.91244433333333
$ws.Range("H10").Value = 113.737333
$ws.Range("I10").Value = 0.543427237729745
$ws.Range("J10").Value = 0.543427237729745
$ws.Range("M10").Value = 11.86694533333333
$ws.Range("N10").Value = 35.600836
$ws.Range("O10").Value = 0.413993613763541
$ws.Range("P10").Value = 0.413993613763541
$ws.Range("Q10").Value = 449.9049043567098
$ws.Range("R10").Value = 4049.144139210388
$ws.Range("S10").Value = 0.224975405965276
$ws.Range("T10").Value = 0.224975405965276

